$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4491677326059005
$ws.Range("C2").Value = 0.2965102684747443
$ws.Range("E2").Value = 0.1539423081504783
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.6457261047731606
$ws.Range("H2").Value = 0.7807710158754659
$ws.Range("K2").Value = 0.2559141323884546
$ws.Range("L2").Value = 0.1902078750966965
$ws.Range("M2").Value = 0.1363245693562938
$ws.Range("O2").Value = 2.843218607440889

$ws.Range("B3").Value = 0.4112595272337103
$ws.Range("C3").Value = 0.2981528487779066
$ws.Range("E3").Value = 0.1548888852018599
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.6540998605119697
$ws.Range("H3").Value = 0.7883942387437131
$ws.Range("K3").Value = 0.2232897513557788
$ws.Range("L3").Value = 0.1875680730085634
$ws.Range("M3").Value = 0.1291374247715247
$ws.Range("O3").Value = 2.876703892095122

$ws.Range("B4").Value = 0.3880274305883518
$ws.Range("C4").Value = 0.2992365000541177
$ws.Range("E4").Value = 0.155538476895881
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.6596516169842843
$ws.Range("H4").Value = 0.7933873130226488
$ws.Range("K4").Value = 0.2031800005366193
$ws.Range("L4").Value = 0.186035923864786
$ws.Range("M4").Value = 0.1247639626971129
$ws.Range("O4").Value = 2.898779997896042

$ws.Range("B5").Value = 0.378571864472633
$ws.Range("C5").Value = 0.2996970296483141
$ws.Range("E5").Value = 0.1558204126056708
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.662017116340472
$ws.Range("H5").Value = 0.7955006381557794
$ws.Range("K5").Value = 0.1949660022575301
$ws.Range("L5").Value = 0.1854339437730914
$ws.Range("M5").Value = 0.1229918102964085
$ws.Range("O5").Value = 2.908157419794641

$ws.Range("B6").Value = 0.3770025027449435
$ws.Range("C6").Value = 0.2997746451605146
$ws.Range("E6").Value = 0.1558682686965813
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.6624161322200379
$ws.Range("H6").Value = 0.7958563034666639
$ws.Range("K6").Value = 0.1936009364676607
$ws.Range("L6").Value = 0.1853353396406874
$ws.Range("M6").Value = 0.1226981580865463
$ws.Range("O6").Value = 2.909737559403176

$ws.Range("B7").Value = 0.3878998611881741
$ws.Range("C7").Value = 0.2992426341944459
$ws.Range("E7").Value = 0.1555422094159837
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.6596831014910549
$ws.Range("H7").Value = 0.7934154957147399
$ws.Range("K7").Value = 0.2030693003154909
$ws.Range("L7").Value = 0.1860277146173459
$ws.Range("M7").Value = 0.1247400218717232
$ws.Range("O7").Value = 2.898904921479101

$ws.Range("B8").Value = 0.4360883307625159
$ws.Range("C8").Value = 0.2970610720525855
$ws.Range("E8").Value = 0.1542545084207454
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.648528205371349
$ws.Range("H8").Value = 0.7833347013509737
$ws.Range("K8").Value = 0.2446818149710452
$ws.Range("L8").Value = 0.1892793076636465
$ws.Range("M8").Value = 0.1338383323804564
$ws.Range("O8").Value = 2.854449699404427

$ws.Range("B9").Value = 0.5309058033286647
$ws.Range("C9").Value = 0.293376784608725
$ws.Range("E9").Value = 0.1522709269501199
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.6299104185462738
$ws.Range("H9").Value = 0.7660423210792686
$ws.Range("K9").Value = 0.3256419666195143
$ws.Range("L9").Value = 0.1963565246354904
$ws.Range("M9").Value = 0.1519879820752728
$ws.Range("O9").Value = 2.779300746289962

$ws.Range("B10").Value = 0.6007344097443763
$ws.Range("C10").Value = 0.2910290054990838
$ws.Range("E10").Value = 0.1511424186074066
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.6182199288327936
$ws.Range("H10").Value = 0.7548432928648694
$ws.Range("K10").Value = 0.384709501897504
$ws.Range("L10").Value = 0.2019802262652775
$ws.Range("M10").Value = 0.1655045000448681
$ws.Range("O10").Value = 2.731418862944892

$ws.Range("B11").Value = 0.6325316513652979
$ws.Range("C11").Value = 0.2900382991092556
$ws.Range("E11").Value = 0.1507001531562437
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.6133338307640557
$ws.Range("H11").Value = 0.7500747091244335
$ws.Range("K11").Value = 0.41148654997329
$ws.Range("L11").Value = 0.2046300462054518
$ws.Range("M11").Value = 0.1716918995376773
$ws.Range("O11").Value = 2.711227329031473

$ws.Range("B12").Value = 0.6445763886262625
$ws.Range("C12").Value = 0.2896742138589445
$ws.Range("E12").Value = 0.1505428790020673
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.6115457601356695
$ws.Range("H12").Value = 0.7483157822657276
$ws.Range("K12").Value = 0.4216124410123427
$ws.Range("L12").Value = 0.2056465622964652
$ws.Range("M12").Value = 0.1740403409366849
$ws.Range("O12").Value = 2.703809987747135

$ws.Range("B13").Value = 0.6419821816416516
$ws.Range("C13").Value = 0.289752134309424
$ws.Range("E13").Value = 0.1505762974071665
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.611928086129069
$ws.Range("H13").Value = 0.7486925162821123
$ws.Range("K13").Value = 0.4194322814059888
$ws.Range("L13").Value = 0.2054270566187597
$ws.Range("M13").Value = 0.1735343239949856
$ws.Range("O13").Value = 2.705397269580857

$ws.Range("B14").Value = 0.6335225071514117
$ws.Range("C14").Value = 0.2900081239093097
$ws.Range("E14").Value = 0.1506870097746216
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.6131854783497843
$ws.Range("H14").Value = 0.7499290627529547
$ws.Range("K14").Value = 0.4123198971239219
$ws.Range("L14").Value = 0.2047134137806239
$ws.Range("M14").Value = 0.1718849997257763
$ws.Range("O14").Value = 2.710612514812595

$ws.Range("B15").Value = 0.6283411858872228
$ws.Range("C15").Value = 0.2901663657811397
$ws.Range("E15").Value = 0.1507561522753349
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.6139637683393175
$ws.Range("H15").Value = 0.750692581140882
$ws.Range("K15").Value = 0.4079615155563374
$ws.Range("L15").Value = 0.204277988820337
$ws.Range("M15").Value = 0.1708754406688442
$ws.Range("O15").Value = 2.713836796314396

$ws.Range("B16").Value = 0.5986569604978911
$ws.Range("C16").Value = 0.2910953024543517
$ws.Range("E16").Value = 0.1511727505250384
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.6185479452897624
$ws.Range("H16").Value = 0.7551614863315308
$ws.Range("K16").Value = 0.3829576336134721
$ws.Range("L16").Value = 0.201808890012444
$ws.Range("M16").Value = 0.1651009039575584
$ws.Range("O16").Value = 2.732770433195327

$ws.Range("B17").Value = 0.5804542614839363
$ws.Range("C17").Value = 0.2916849450180052
$ws.Range("E17").Value = 0.1514465146800923
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.6214708777298483
$ws.Range("H17").Value = 0.757986464653257
$ws.Range("K17").Value = 0.3675942879428931
$ws.Range("L17").Value = 0.2003175736286096
$ws.Range("M17").Value = 0.1615682076187142
$ws.Range("O17").Value = 2.744792951764566

$ws.Range("B18").Value = 0.5699875949793238
$ws.Range("C18").Value = 0.2920313711838318
$ws.Range("E18").Value = 0.1516106700248194
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.6231927253789138
$ws.Range("H18").Value = 0.7596419956566436
$ws.Range("K18").Value = 0.3587489741550201
$ws.Range("L18").Value = 0.1994684317204616
$ws.Range("M18").Value = 0.1595399428330495
$ws.Range("O18").Value = 2.751857655966077

$ws.Range("B19").Value = 0.5664443079108992
$ws.Range("C19").Value = 0.2921499166947612
$ws.Range("E19").Value = 0.1516674004706076
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.6237826935837276
$ws.Range("H19").Value = 0.7602078000343226
$ws.Range("K19").Value = 0.3557526232853547
$ws.Range("L19").Value = 0.1991824104727442
$ws.Range("M19").Value = 0.1588538385574978
$ws.Range("O19").Value = 2.754275346577117

$ws.Range("B20").Value = 0.5823916603497423
$ws.Range("C20").Value = 0.291621423489012
$ws.Range("E20").Value = 0.1514166794295999
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.6211555184939499
$ws.Range("H20").Value = 0.7576825660524449
$ws.Range("K20").Value = 0.3692306496312767
$ws.Range("L20").Value = 0.2004754348544679
$ws.Range("M20").Value = 0.16194389265479
$ws.Range("O20").Value = 2.743497643258763

$ws.Range("B21").Value = 0.6360072212817158
$ws.Range("C21").Value = 0.289932633363442
$ws.Range("E21").Value = 0.1506542141578997
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.6128144634158303
$ws.Range("H21").Value = 0.7495645883277575
$ws.Range("K21").Value = 0.4144093615114741
$ws.Range("L21").Value = 0.2049226734547176
$ws.Range("M21").Value = 0.172369300703501
$ws.Range("O21").Value = 2.709074461832728

$ws.Range("B22").Value = 0.6710699529677413
$ws.Range("C22").Value = 0.2888934378908701
$ws.Range("E22").Value = 0.1502153549968099
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.6077256182077377
$ws.Range("H22").Value = 0.744531979975541
$ws.Range("K22").Value = 0.4438543995978819
$ws.Range("L22").Value = 0.2079054431591771
$ws.Range("M22").Value = 0.1792143557271828
$ws.Range("O22").Value = 2.687910319433456

$ws.Range("B23").Value = 0.652354566415454
$ws.Range("C23").Value = 0.2894421861079053
$ws.Range("E23").Value = 0.1504441493519586
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.6104084365150513
$ws.Range("H23").Value = 0.7471930128518594
$ws.Range("K23").Value = 0.4281467208511742
$ws.Range("L23").Value = 0.2063065334171057
$ws.Range("M23").Value = 0.1755581939932895
$ws.Range("O23").Value = 2.699083985366869

$ws.Range("B24").Value = 0.5815157675359615
$ws.Range("C24").Value = 0.2916501184173512
$ws.Range("E24").Value = 0.151430146873734
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.6212979634065121
$ws.Range("H24").Value = 0.7578198607568538
$ws.Range("K24").Value = 0.3684908900929145
$ws.Range("L24").Value = 0.2004040401351688
$ws.Range("M24").Value = 0.1617740369369329
$ws.Range("O24").Value = 2.744082776287811

$ws.Range("B25").Value = 0.5052238972379257
$ws.Range("C25").Value = 0.2943102182448456
$ws.Range("E25").Value = 0.1527496954137533
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.6345980877416295
$ws.Range("H25").Value = 0.7704556834763565
$ws.Range("K25").Value = 0.303811160104317
$ws.Range("L25").Value = 0.194367248016512
$ws.Range("M25").Value = 0.1470456520947252
$ws.Range("O25").Value = 2.798343128562706
